# feat: add 2022-Q4 data
#
# This script reshapes the workbook so that:
#   - "总计" summary sheet gets a new "2022-Q4" row inserted (and the older
#     rows shift down to make room for it).
#   - the sheet that used to hold the "2022-Q3" per-fund detail becomes the
#     new "2022-Q4" detail sheet (with brand-new fund-holding data).
#   - a new sheet named "2022-Q3" is inserted after it, carrying the detail
#     data that used to live in the (now renamed) "2022-Q4" sheet.
#   - the "2021-Q4" detail sheet is left with the same data, shifted to the
#     end (position 4).

$wb = $excel.ActiveWorkbook

# Helper: assign a value as TEXT (matching the workbook's convention of
# storing these columns as inline/shared strings, even when the text looks
# like a number) without leaving a visible "quote prefix" style behind.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ------------------------------------------------------------------
# Step 1: duplicate the current "2022-Q3" detail sheet (position 2).
# The duplicate will become the new "2022-Q3" sheet (position 3),
# preserving the original fund-holding rows (005997 / 002388).
# The original sheet object (position 2) will be turned into the new
# "2022-Q4" sheet with fresh data.
# ------------------------------------------------------------------
$sheetQ3 = $wb.Worksheets.Item(2)
$sheetQ3.Copy($null, $sheetQ3)
$newSheetQ3 = $wb.Worksheets.Item(3)

# Rename the original first (frees up the "2022-Q3" name), then rename the
# duplicate to take over that freed-up name.
$sheetQ3.Name = "2022-Q4"
$newSheetQ3.Name = "2022-Q3"

# ------------------------------------------------------------------
# Step 2: replace the data on the (renamed) "2022-Q4" sheet with the new
# single-fund holding, and drop the now-unused second data row.
# ------------------------------------------------------------------
Set-TextValue $sheetQ3.Range("B2") "506009"
Set-TextValue $sheetQ3.Range("C2") "国泰科创板两年定期开放混合"
Set-TextValue $sheetQ3.Range("D2") "2.17"
Set-TextValue $sheetQ3.Range("E2") "87.42"
Set-TextValue $sheetQ3.Range("F2") "2.95"
Set-TextValue $sheetQ3.Range("G2") "0.0640"
$sheetQ3.Range("H2").Value = 6
$sheetQ3.Rows.Item(3).Delete()

# ------------------------------------------------------------------
# Step 3: update the "总计" (summary) sheet. A new "2022-Q4" row is
# inserted right after the header, pushing the existing quarters down by
# one row.
# ------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item(1)

# Grab the old row2/row3 values before we overwrite anything.
$oldRow2B = $sheetTotal.Range("B2").Value()
$oldRow2C = $sheetTotal.Range("C2").Value()
$oldRow2D = $sheetTotal.Range("D2").Value()
$oldRow3B = $sheetTotal.Range("B3").Value()
$oldRow3C = $sheetTotal.Range("C3").Value()
$oldRow3D = $sheetTotal.Range("D3").Value()

# Build the new row4 (style copied from row3 so A4 keeps the same look).
$sheetTotal.Range("A3").Copy($sheetTotal.Range("A4"))
$sheetTotal.Range("A4").Value = 2
$sheetTotal.Range("B4").Value = $oldRow3B
$sheetTotal.Range("C4").Value = $oldRow3C
$sheetTotal.Range("D4").Value = $oldRow3D

# Old row2 data moves down to row3.
$sheetTotal.Range("B3").Value = $oldRow2B
$sheetTotal.Range("C3").Value = $oldRow2C
$sheetTotal.Range("D3").Value = $oldRow2D

# New "2022-Q4" data goes into row2.
$sheetTotal.Range("B2").Value = "2022-Q4"
$sheetTotal.Range("C2").Value = 1
$sheetTotal.Range("D2").Value = 0.06
